# Update "want to go" (想去人数) counts (column F) across all sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 851
$ws1.Range("F3").Value = 1763
$ws1.Range("F4").Value = 52
$ws1.Range("F5").Value = 534
$ws1.Range("F6").Value = 2147
$ws1.Range("F7").Value = 1373
$ws1.Range("F8").Value = 2088
$ws1.Range("F9").Value = 970
$ws1.Range("F11").Value = 2406
$ws1.Range("F12").Value = 661
$ws1.Range("F14").Value = 3935
$ws1.Range("F16").Value = 368
$ws1.Range("F17").Value = 3028
$ws1.Range("F18").Value = 812
$ws1.Range("F19").Value = 144
$ws1.Range("F21").Value = 127
$ws1.Range("F22").Value = 2059
$ws1.Range("F23").Value = 1181
$ws1.Range("F24").Value = 1899
$ws1.Range("F25").Value = 389
$ws1.Range("F26").Value = 208
$ws1.Range("F27").Value = 14
$ws1.Range("F28").Value = 8454
$ws1.Range("F29").Value = 5703
$ws1.Range("F30").Value = 354
$ws1.Range("F31").Value = 177
$ws1.Range("F32").Value = 752
$ws1.Range("F33").Value = 768
$ws1.Range("F34").Value = 3476
$ws1.Range("F36").Value = 946
$ws1.Range("F37").Value = 390
$ws1.Range("F38").Value = 37
$ws1.Range("F41").Value = 4630
$ws1.Range("F42").Value = 1
$ws1.Range("F43").Value = 856
$ws1.Range("F44").Value = 79
$ws1.Range("F45").Value = 400

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 182
$ws2.Range("F19").Value = 66
$ws2.Range("F25").Value = 7
$ws2.Range("F26").Value = 29

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8295
$ws3.Range("F4").Value = 1292

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 851
$ws4.Range("F4").Value = 1292
$ws4.Range("F6").Value = 1763
$ws4.Range("F7").Value = 52
$ws4.Range("F8").Value = 534
$ws4.Range("F9").Value = 1373
$ws4.Range("F10").Value = 2088
$ws4.Range("F11").Value = 970
$ws4.Range("F15").Value = 3935
$ws4.Range("F16").Value = 368
$ws4.Range("F17").Value = 3028
$ws4.Range("F18").Value = 812
$ws4.Range("F19").Value = 144
$ws4.Range("F21").Value = 2059
$ws4.Range("F27").Value = 1899
$ws4.Range("F29").Value = 208
$ws4.Range("F30").Value = 14
$ws4.Range("F31").Value = 8454
$ws4.Range("F32").Value = 5703
$ws4.Range("F33").Value = 66
$ws4.Range("F34").Value = 354
$ws4.Range("F35").Value = 177
$ws4.Range("F36").Value = 752
$ws4.Range("F37").Value = 768
$ws4.Range("F39").Value = 946
$ws4.Range("F40").Value = 390
$ws4.Range("F43").Value = 4630
$ws4.Range("F44").Value = 856
$ws4.Range("F45").Value = 400
$ws4.Range("F47").Value = 7
$ws4.Range("F48").Value = 29
